# Weekly price-data update: a new "Albahaca" price record for
# Vega Modelo de Temuco is inserted at row 309 (newest entry on top),
# shifting the existing history rows 309-376 down to 310-377.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(309).Insert()

$ws.Range("A309").Value = 10
$ws.Range("B309").Value = "Vega Modelo de Temuco"
$ws.Range("C309").Value = "La Araucanía"
$ws.Range("D309").Value = 45015
$ws.Range("E309").Value = 9
$ws.Range("F309").Value = 100112052
$ws.Range("G309").Value = "Albahaca"
$ws.Range("H309").Value = "Sin especificar"
$ws.Range("I309").Value = "Primera"
$ws.Range("J309").Value = 80
$ws.Range("K309").Value = 5000
$ws.Range("L309").Value = 6000
$ws.Range("M309").Value = 5500
$ws.Range("N309").Value = "$/paquete"
$ws.Range("O309").Value = "Región de La Araucanía"
$ws.Range("P309").Value = 5500
$ws.Range("Q309").Value = 1
$ws.Range("R309").Value = "Hortaliza"
